# Doing Updates for Financials
# Update the Balance Sheet figures (column D = most recent period) on the PRLX sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Other Current Assets
$ws.Range("D45").Value = 100

# Goodwill
$ws.Range("D49").Value = 1500

# Other Assets
$ws.Range("D52").Value = 200

# Total Assets
$ws.Range("D54").Value = 1800

# Accounts Payable
$ws.Range("D57").Value = 3000

# Other Current Liabilities
$ws.Range("D59").Value = 600

# Total Current Liabilities
$ws.Range("D60").Value = 4400

# Long Term Debt
$ws.Range("D61").Value = 19900

# Other Liabilities - now populated for the latest period, older periods become NA
$ws.Range("D62").Value = 4600
$ws.Range("E62").Value = "NA"
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "NA"
$ws.Range("H62").Value = "NA"
$ws.Range("I62").Value = "NA"
$ws.Range("J62").Value = "NA"

# Total Liabilities
$ws.Range("D66").Value = 28900

# Net Income
$ws.Range("D72").Value = -33600

# Total Stockholder Equity
$ws.Range("D76").Value = -27200
